$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare new text content (order matters for shared-string indexing) ---

$filesTabLabel = "FilesTab"

$casesQuery = @"
MATCH (c:case)
 MATCH (c)-[:of_arm]->(a:arm)-[:of_trial]->(ct:clinical_trial)
 MATCH (f:file)-[*]->(c)
   WHERE c.race = "WHITE"
RETURN DISTINCT
    c.case_id AS ``Case ID``,
     ct.clinical_trial_designation AS ``Trial Code``,
     a.arm_id AS Arm,
      a.arm_drug AS ``Arm Treatment``,
c.disease AS Diagnosis,
  c.gender AS Gender,
    c.race AS Race,
    c.ethnicity AS Ethnicity
"@

$statsQuery = @"
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
WITH f,a,ct,c
        WHERE c.race = "WHITE"
RETURN
    COUNT(DISTINCT ct.clinical_trial_designation) AS Trials,
    COUNT(DISTINCT c.case_id) AS Cases,
    COUNT(DISTINCT f) AS Files
"@

$filesQuery = @"
MATCH (f:file)
OPTIONAL MATCH (f)-[*]->(a:arm)-[:of_trial]->(ct:clinical_trial)
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (f)-->(parent)
WITH f,a,ct,c,parent
  WHERE c.race = "WHITE"
WITH
    f, parent, c, a, ct,
    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
    toInteger(floor(log(f.file_size)/log(1024))) as i,
    2 as precision
WITH
    f, parent, c, a, ct,
    f.file_size /(1024^i) AS value,
    10^precision AS factor,
    units[i] as unit
WITH
    f, parent, c, a, ct, unit,
    round(factor * value)/factor AS size
RETURN DISTINCT
    f.file_name AS ``File Name``,
    head(labels(parent)) as Association,
    f.file_description AS Description,
    f.file_format AS ``File Format``,
    CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    ct.clinical_trial_designation AS ``Trial Code``,
    a.arm_id AS Arm,
    c.case_id AS ``Case ID``
"@

# --- Row 3 label first, then update Row 2 query cells, then Row 3 query cells ---
# (keeps shared-string insertion order aligned with target workbook)
$ws.Range("A3").Value = $filesTabLabel
$ws.Range("B2").Value = $casesQuery
$ws.Range("C2").Value = $statsQuery
$ws.Range("B3").Value = $filesQuery
$ws.Range("C3").Value = $statsQuery
$ws.Range("D3").Value = $ws.Range("D2").Value()
$ws.Range("E3").Value = $ws.Range("E2").Value()

# Apply the same wrap-text formatting used by B2/C2 to B3/C3
$ws.Range("B3").WrapText = $true
$ws.Range("C3").WrapText = $true

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 195
$ws.Rows.Item(3).RowHeight = 409.5

# --- Sheet view settings: zoom 70%, scroll to row 3, select D3 ---
$window = $excel.ActiveWindow
$window.ScrollRow = 3
$window.ScrollColumn = 1
$window.Zoom = 70
$ws.Range("D3").Select()
